$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Delete()

for ($col = 2; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $cell.Value2.ToString() + ".global"
}
